$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet/workbook title to reflect new "through" date
$ws.Name = "Through 2022-03-06"

# Update the row label for March
$ws.Range("A4").Value = "March (through 03-06)"

# Update March row (row 4) values
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 16
$ws.Range("I4").Value = 33

# Update Total row (row 5) values
$ws.Range("B5").Value = 42
$ws.Range("C5").Value = 94
$ws.Range("D5").Value = 140
$ws.Range("E5").Value = 149
$ws.Range("F5").Value = 85
$ws.Range("G5").Value = 155
$ws.Range("H5").Value = 358
$ws.Range("I5").Value = 334
